$wb = $excel.ActiveWorkbook

$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo.Range("D7").Value = 442.27
$wsVentasGrupo.Range("L7").Value = 86.56999999999999
$wsVentasGrupo.Range("M7").Value = 612.86
$wsVentasGrupo.Range("O7").Value = 1070.53

$wsVentasGrupo.Range("H15").Value = 603
$wsVentasGrupo.Range("I15").Value = 28.8
$wsVentasGrupo.Range("Q15").Value = 97.11

$wsVentasGrupo.Range("D22").Value = "1 de 20"
$wsVentasGrupo.Range("H22").Value = "1 de 20"
$wsVentasGrupo.Range("I22").Value = "1 de 20"
$wsVentasGrupo.Range("L22").Value = "1 de 20"
$wsVentasGrupo.Range("M22").Value = "1 de 20"
$wsVentasGrupo.Range("O22").Value = "1 de 20"
$wsVentasGrupo.Range("Q22").Value = "1 de 20"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual.Range("F7").Value = 2212.23
$wsVentaMensual.Range("F15").Value = 728.91
$wsVentaMensual.Range("F22").Value = 2941.14

$wsVentaMensual.Columns.Item(6).ColumnWidth = 12.166666666666666
